$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, pushing existing rows 58..148 down to 59..149
$ws.Rows("58").Insert()

# Populate the newly inserted row 58 with the new record
$ws.Range("A58").Value = 4
$ws.Range("B58").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C58").Value = "Los Lagos"
$ws.Range("D58").Value = 44540
$ws.Range("E58").Value = 10
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100102
$ws.Range("H58").Value = "Cítricos"
$ws.Range("I58").Value = 100102004
$ws.Range("J58").Value = "Mandarina"
$ws.Range("K58").Value = "Murcott"
$ws.Range("L58").Value = "Segunda"
$ws.Range("M58").Value = 900
$ws.Range("N58").Value = 13000
$ws.Range("O58").Value = 13000
$ws.Range("P58").Value = 13000
$ws.Range("Q58").Value = "$/caja 18 kilos"
$ws.Range("R58").Value = "Región Metropolitana"
$ws.Range("S58").Value = 722
$ws.Range("T58").Value = 18
